$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 11, column B (the 4th rule row) held the text "R40"; it must become
# the text "1" - a genuine (shared-string) text value, not a number, and
# without disturbing B11's existing cell style (s="23").
#
# A plain Range.Value assignment of "1" is auto-coerced to a number by the
# COM layer, and forcing text via a leading apostrophe / NumberFormat="@"
# both mint a brand-new (quote-prefixed / text) cell style, which would
# leave an extra, unused entry behind in styles.xml. Instead, build the
# text value with a self-referential formula (always textual, t="str"),
# copy *only the value* onto B11 via PasteSpecial so its style is left
# completely untouched, then remove the scratch cell again.
$helper = $ws.Range("Z100")
$helper.Formula = "=""1"""
$helper.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$helper.Delete(-4162)                 # xlShiftUp - remove the scratch cell entirely
